$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# Update price in C2 (97 -> 1.49)
$ws.Range("C2").Value = 1.49

# Add new row of data for "Captain America t shirt"
$ws.Range("B3").Value = "Captain America t shirt"
$ws.Range("C3").Value = 2.45
$ws.Range("D3").Value = 300
$ws.Range("E3").Value = 280
$ws.Range("F3").Value = 3.91
$ws.Range("G3").Value = 150

# Move the active selection to G4, matching the author's final cursor position
$ws.Activate()
$ws.Range("G4").Select()
